$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TestCaseNameDesc (column D) text updates ---
$ws.Range("D2").Value = "Validate pet response"
$ws.Range("D3").Value = "Get by Id"
$ws.Range("D4").Value = "Get api testing"
$ws.Range("D5").Value = "Post api testing with kafka aggregated "
$ws.Range("D6").Value = "Message verification"
$ws.Range("D7").Value = "Kafka_aggregated"
$ws.Range("D8").Value = "Post proto buff message"
$ws.Range("D9").Value = "Validate Protobuff message"

# --- StepInfo (column E) text updates ---
$ws.Range("E6").Value = "contains json based information"
$ws.Range("E7").Value = "contains json based information"
$ws.Range("E9").Value = "contains proto buff based information"

# --- Column width adjustments (D widened, E narrowed) to fit the new text ---
$ws.Columns("D").ColumnWidth = 33.33
$ws.Columns("E").ColumnWidth = 31.5

# --- Selection / scroll position moves from P8 to D7 ---
$ws.Range("D7").Select()
